{"js": "const replacements = [\n  [\"80-10=70\", \"67-47=20\"],\n  [\"45-10=35\", \"89-71=18\"],\n  [\"62-31=31\", \"22+11=33\"],\n  [\"40+21=61\", \"29+16=45\"],\n  [\"35+10=45\", \"39+23=62\"],\n  [\"75-24=51\", \"86-1=85\"],\n  [\"42+18=60\", \"92-46=46\"],\n  [\"56+12=68\", \"83-34=49\"],\n  [\"34+17=51\", \"20+32=52\"],\n  [\"97-10=87\", \"36-1=35\"],\n  [\"47-21=26\", \"8+74=82\"],\n  [\"76+22=98\", \"9+65=74\"],\n  [\"69+15=84\", \"38+26=64\"],\n  [\"35+63=98\", \"16+61=77\"],\n  [\"44+47=91\", \"57+35=92\"],\n  [\"79-68=11\", \"42-11=31\"],\n  [\"97-42=55\", \"66+0=66\"],\n  [\"5+70=75\", \"76+1=77\"],\n  [\"7+51=58\", \"86-85=1\"],\n  [\"23+76=99\", \"6+91=97\"],\n  [\"29+57=86\", \"5+86=91\"],\n  [\"7-4=3\", \"69-54=15\"],\n  [\"49+47=96\", \"98-13=85\"],\n  [\"69-34=35\", \"2+13=15\"],\n  [\"10+36=46\", \"86+13=99\"],\n  [\"1+95=96\", \"20+76=96\"],\n  [\"8+28=36\", \"76-6=70\"],\n  [\"78-37=41\", \"4+75=79\"],\n  [\"91+4=95\", \"97-16=81\"],\n  [\"33+33=66\", \"82-3=79\"],\n  [\"13+85=98\", \"10+2=12\"],\n  [\"77-52=25\", \"87-74=13\"],\n  [\"56+20=76\", \"96-1=95\"],\n  [\"32+47=79\", \"49+5=54\"],\n  [\"89-21=68\", \"14-10=4\"],\n  [\"68+10=78\", \"4+50=54\"],\n  [\"60+19=79\", \"72-62=10\"],\n  [\"31-23=8\", \"93-57=36\"],\n  [\"82-24=58\", \"20+37=57\"],\n  [\"40-35=5\", \"74-57=17\"],\n  [\"50-16=34\", \"89-78=11\"],\n  [\"8+77=85\", \"94-71=23\"],\n  [\"95-46=49\", \"35-4=31\"],\n  [\"55+27=82\", \"82+13=95\"],\n  [\"41+11=52\", \"69-11=58\"],\n  [\"87-67=20\", \"60-48=12\"],\n  [\"17+30=47\", \"68-67=1\"],\n  [\"75-1=74\", \"31+45=76\"],\n  [\"42+44=86\", \"77-25=52\"],\n  [\"11+15=26\", \"2+58=60\"],\n  [\"83+9=92\", \"34+61=95\"],\n  [\"69-52=17\", \"59+20=79\"],\n  [\"82-70=12\", \"31-15=16\"],\n  [\"89-33=56\", \"1+47=48\"],\n  [\"89-74=15\", \"64-43=21\"],\n  [\"26-22=4\", \"68-68=0\"],\n  [\"62+1=63\", \"11-4=7\"],\n  [\"71-9=62\", \"27+63=90\"],\n  [\"33+20=53\", \"53-6=47\"],\n  [\"72-23=49\", \"76-16=60\"],\n  [\"7+90=97\", \"78+6=84\"],\n  [\"21+60=81\", \"81-4=77\"],\n  [\"34-25=9\", \"86-9=77\"],\n  [\"9+17=26\", \"41+38=79\"],\n  [\"36+37=73\", \"66-23=43\"],\n  [\"25+27=52\", \"67-64=3\"],\n  [\"22-8=14\", \"40+58=98\"],\n  [\"55-27=28\", \"38-2=36\"],\n  [\"17+53=70\", \"1+94=95\"],\n  [\"10+54=64\", \"4+54=58\"],\n  [\"37-29=8\", \"37-23=14\"],\n  [\"63+17=80\", \"87-85=2\"],\n  [\"36-8=28\", \"89-9=80\"],\n  [\"27-14=13\", \"94+0=94\"],\n  [\"21+8=29\", \"39+42=81\"],\n  [\"25+63=88\", \"13+71=84\"],\n  [\"13+23=36\", \"93-92=1\"],\n  [\"61-23=38\", \"74-52=22\"],\n  [\"72-3=69\", \"40+40=80\"],\n  [\"30+36=66\", \"98-48=50\"],\n  [\"40+6=46\", \"99-89=10\"],\n  [\"65-4=61\", \"9+48=57\"],\n  [\"55-45=10\", \"57+15=72\"],\n  [\"59-47=12\", \"71-64=7\"],\n  [\"11+67=78\", \"71-62=9\"],\n  [\"27+42=69\", \"39+52=91\"],\n  [\"92-33=59\", \"91-32=59\"],\n  [\"86-57=29\", \"98-11=87\"],\n  [\"30+29=59\", \"93-50=43\"],\n  [\"62-62=0\", \"43-5=38\"],\n  [\"28+65=93\", \"60-34=26\"],\n  [\"26+56=82\", \"7+52=59\"],\n  [\"38+42=80\", \"46+10=56\"],\n  [\"92-51=41\", \"30-12=18\"],\n  [\"10+74=84\", \"66+17=83\"],\n  [\"47-2=45\", \"16+53=69\"],\n  [\"64+15=79\", \"31-1=30\"],\n  [\"55-11=44\", \"15+12=27\"],\n  [\"67-17=50\", \"34+38=72\"],\n  [\"83-72=11\", \"35+3=38\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@('80-10=70', '67-47=20')\n    ,@('45-10=35', '89-71=18')\n    ,@('62-31=31', '22+11=33')\n    ,@('40+21=61', '29+16=45')\n    ,@('35+10=45', '39+23=62')\n    ,@('75-24=51', '86-1=85')\n    ,@('42+18=60', '92-46=46')\n    ,@('56+12=68', '83-34=49')\n    ,@('34+17=51', '20+32=52')\n    ,@('97-10=87', '36-1=35')\n    ,@('47-21=26', '8+74=82')\n    ,@('76+22=98', '9+65=74')\n    ,@('69+15=84', '38+26=64')\n    ,@('35+63=98', '16+61=77')\n    ,@('44+47=91', '57+35=92')\n    ,@('79-68=11', '42-11=31')\n    ,@('97-42=55', '66+0=66')\n    ,@('5+70=75', '76+1=77')\n    ,@('7+51=58', '86-85=1')\n    ,@('23+76=99', '6+91=97')\n    ,@('29+57=86', '5+86=91')\n    ,@('7-4=3', '69-54=15')\n    ,@('49+47=96', '98-13=85')\n    ,@('69-34=35', '2+13=15')\n    ,@('10+36=46', '86+13=99')\n    ,@('1+95=96', '20+76=96')\n    ,@('8+28=36', '76-6=70')\n    ,@('78-37=41', '4+75=79')\n    ,@('91+4=95', '97-16=81')\n    ,@('33+33=66', '82-3=79')\n    ,@('13+85=98', '10+2=12')\n    ,@('77-52=25', '87-74=13')\n    ,@('56+20=76', '96-1=95')\n    ,@('32+47=79', '49+5=54')\n    ,@('89-21=68', '14-10=4')\n    ,@('68+10=78', '4+50=54')\n    ,@('60+19=79', '72-62=10')\n    ,@('31-23=8', '93-57=36')\n    ,@('82-24=58', '20+37=57')\n    ,@('40-35=5', '74-57=17')\n    ,@('50-16=34', '89-78=11')\n    ,@('8+77=85', '94-71=23')\n    ,@('95-46=49', '35-4=31')\n    ,@('55+27=82', '82+13=95')\n    ,@('41+11=52', '69-11=58')\n    ,@('87-67=20', '60-48=12')\n    ,@('17+30=47', '68-67=1')\n    ,@('75-1=74', '31+45=76')\n    ,@('42+44=86', '77-25=52')\n    ,@('11+15=26', '2+58=60')\n    ,@('83+9=92', '34+61=95')\n    ,@('69-52=17', '59+20=79')\n    ,@('82-70=12', '31-15=16')\n    ,@('89-33=56', '1+47=48')\n    ,@('89-74=15', '64-43=21')\n    ,@('26-22=4', '68-68=0')\n    ,@('62+1=63', '11-4=7')\n    ,@('71-9=62', '27+63=90')\n    ,@('33+20=53', '53-6=47')\n    ,@('72-23=49', '76-16=60')\n    ,@('7+90=97', '78+6=84')\n    ,@('21+60=81', '81-4=77')\n    ,@('34-25=9', '86-9=77')\n    ,@('9+17=26', '41+38=79')\n    ,@('36+37=73', '66-23=43')\n    ,@('25+27=52', '67-64=3')\n    ,@('22-8=14', '40+58=98')\n    ,@('55-27=28', '38-2=36')\n    ,@('17+53=70', '1+94=95')\n    ,@('10+54=64', '4+54=58')\n    ,@('37-29=8', '37-23=14')\n    ,@('63+17=80', '87-85=2')\n    ,@('36-8=28', '89-9=80')\n    ,@('27-14=13', '94+0=94')\n    ,@('21+8=29', '39+42=81')\n    ,@('25+63=88', '13+71=84')\n    ,@('13+23=36', '93-92=1')\n    ,@('61-23=38', '74-52=22')\n    ,@('72-3=69', '40+40=80')\n    ,@('30+36=66', '98-48=50')\n    ,@('40+6=46', '99-89=10')\n    ,@('65-4=61', '9+48=57')\n    ,@('55-45=10', '57+15=72')\n    ,@('59-47=12', '71-64=7')\n    ,@('11+67=78', '71-62=9')\n    ,@('27+42=69', '39+52=91')\n    ,@('92-33=59', '91-32=59')\n    ,@('86-57=29', '98-11=87')\n    ,@('30+29=59', '93-50=43')\n    ,@('62-62=0', '43-5=38')\n    ,@('28+65=93', '60-34=26')\n    ,@('26+56=82', '7+52=59')\n    ,@('38+42=80', '46+10=56')\n    ,@('92-51=41', '30-12=18')\n    ,@('10+74=84', '66+17=83')\n    ,@('47-2=45', '16+53=69')\n    ,@('64+15=79', '31-1=30')\n    ,@('55-11=44', '15+12=27')\n    ,@('67-17=50', '34+38=72')\n    ,@('83-72=11', '35+3=38')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $result = $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $result) {\n        throw \"Replacement failed for: $oldText\"\n    }\n}"}
